# Append new experiment records (2020-03-18, RF mapping + reduced-dimension evolution comparison)
# for the second monkey dataset, per commit: "Support sorting exp record to 2 monkey's different dests"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A179").Value = "Beto-18032020-001"
$ws.Range("B179").Value = "200318_Beto_rfMapper_basic"
$ws.Range("C179").Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-18-Beto"
$ws.Range("D179").Value = "001 RFMapping 11:`n-8:8:8 `nCarlos' huge image test`nCompleted"
$ws.Range("E179").Value = "ReducDimen_Evol"

$ws.Range("A180").Value = "Beto-18032020-002"
$ws.Range("B180").Value = "200318_Beto_rfMapper_basic(1)"
$ws.Range("C180").Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-18-Beto"
$ws.Range("D180").Value = "002 RFMapping starts 11:35`n-8:2:8`nCompleted"
$ws.Range("E180").Value = "ReducDimen_Evol"

$ws.Range("A181").Value = "Beto-18032020-003"
$ws.Range("B181").Value = "200318_Beto_rfMapper_basic(2)"
$ws.Range("C181").Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-18-Beto"
$ws.Range("D181").Value = "003 RFMapping starts 11:49`n-4:1:4`n8 mins similar RF position`nCompleted"
$ws.Range("E181").Value = "ReducDimen_Evol"

$ws.Range("A182").Value = "Beto-18032020-004"
$ws.Range("B182").Value = "200318_Beto_generate_integrated(1)"
$ws.Range("C182").Value = "N:\Stimuli\2019-12-Evolutions\2020-03-18-Beto-01\2020-03-18-12-07-15"
$ws.Range("D182").Value = "004 Generate Integrated`n24 [-1.5 1.5 ] 4 1 CMAES`nTrying out basic optimizer, see if it can evolve!`nSeems work, Evolving very fast.`nSeems plateau really fast as well around 10 gens. Only 6 mins! Wow. Maybe 20 blocks`nWhy Stuck so early? Gets some local maxima? `nCompleted"
$ws.Range("E182").Value = "ReducDimen_Evol"

$ws.Range("A183").Value = "Beto-18032020-005"
$ws.Range("B183").Value = "200318_Beto_generate_integrated(2)"
$ws.Range("C183").Value = "N:\Stimuli\2019-12-Evolutions\2020-03-18-Beto-02\2020-03-18-12-21-38"
$ws.Range("D183").Value = "005 Generate Integrated`n24 [-1.5 1.5 ] 4 1 ZOHA_Sphere_lr_euclid`n24 [-1.5 1.5 ] 4 1 ZOHA_Sphere_lr_euclid_RD`nTest the reduced Dimension Comparison.`nUse the inverse decay exploration range. lr=1.5`nInverse decay's initial learning rate is too large, should make it smaller. And learning rate * exploration exceeds pi/2 it's insane.`n(Gen 5 explor is 58.3 deg, step is 0.88 or 0.64)`nGen 12 42 degs, step 0.474, 0.314, this is good.`nGap exist! But quite small. (only 14 mins now! )`nPlateau around 15 gens`nCompleted"
$ws.Range("E183").Value = "ReducDimen_Evol"

$ws.Range("A184").Value = "Beto-18032020-006"
$ws.Range("B184").Value = "200318_Beto_generate_integrated(3)"
$ws.Range("C184").Value = "N:\Stimuli\2019-12-Evolutions\2020-03-18-Beto-03\2020-03-18-12-51-20"
$ws.Range("D184").Value = "006 Generate Integrated`n24 [-1.5 1.5 ] 4 1 ZOHA_Sphere_lr_euclid`n24 [-1.5 1.5 ] 4 1 ZOHA_Sphere_lr_euclid_RD`nRedo the reduced Dimension Comparison.`nUse the inverse decay exploration range. Lr=1.2 Decrease this to avoid overshoot. And this makes Sphere lr euclid grow just as fast as CMAES! No delay!`ngen 11, 45 deg 0.321.`nInterestingly, the reduced starts to grow at around gen 12. Making gap smaller.`nFinished in 30 mins`nCompleted"
$ws.Range("E184").Value = "ReducDimen_Evol"

$ws.Range("A185").Value = "Beto-18032020-007"
$ws.Range("B185").Value = "200318_Beto_generate_integrated(4)"
$ws.Range("C185").Value = "N:\Stimuli\2019-12-Evolutions\2020-03-18-Beto-04\2020-03-18-13-30-57"
$ws.Range("D185").Value = "007 Generate Integrated starts 13:30`n64 [-1.5 -2.5 ] 5 1 ZOHA_Sphere_lr_euclid`n64 [-1.5 -2.5 ] 5 1 ZOHA_Sphere_lr_euclid_RD`nTest the Reduced Dimension for V4 channel on Beto.`nZOHA full saturates pretty fast!`nSeems there is still a gap and it's not small……`nWow, gap gets closed…… as expected.`nReally the same as expected! It closes the gap finally.`ntakes 35 min to get to 25 blocks! Good job.`nCompleted"
$ws.Range("E185").Value = "ReducDimen_Evol"

$ws.Range("A186").Value = "Beto-18032020-008"
$ws.Range("B186").Value = "200318_Beto_generate_integrated(5)"
$ws.Range("C186").Value = "N:\Stimuli\2019-12-Evolutions\2020-03-18-Beto-05\2020-03-18-14-10-20"
$ws.Range("D186").Value = "008 Generate Integrated starts 14:10`n64 [-1.5 -2.5 ] 5 1 ZOHA_Sphere_lr_euclid_RD`n64 [-1.5 -2.5 ] 5 1 ZOHA_Sphere_lr_euclid_RD`nTest the Reduced Dimension for V4 channel on Beto.`nSee how large is the trial variability.`nThe first one saturates pretty fast as well.`nthe first one gets even higher than full, seemingly! OMG.`nAt gen 10 the step size is around 0.376!`nOh get back to normal… seems higher than Full evolution is just fluctuation.`nGen 13 exploration aaround 41.8.`nAround 18 gens the 1st evolution grows again! Find something new to add to the image.`nGets to ~ 30 gens……`nAdd 100mL water to him!`nCompleted"
$ws.Range("E186").Value = "ReducDimen_Evol"

